# Applies the "Updated symbol list" diff to the crypto price sheet.
# Rewrites B/C/D/E cells (Coin, Link, Price, Volume(1h)) for the rows
# that changed, keeping every value as literal text (matching the
# workbook's inline-string cells) rather than letting Excel coerce
# numeric-looking strings (e.g. "243.73", "0.0002000") into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "243.73"
Set-TextValue $ws "D3" "24.97"
Set-TextValue $ws "B4" "HuobiToken"
Set-TextValue $ws "C4" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws "D4" "5.164"
Set-TextValue $ws "E4" "3HuobiTokenHT"
Set-TextValue $ws "B5" "Cronos"
Set-TextValue $ws "C5" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D5" "0.05715"
Set-TextValue $ws "E5" "4CronosCRO"
Set-TextValue $ws "B6" "KuCoinToken"
Set-TextValue $ws "C6" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws "D6" "6.472"
Set-TextValue $ws "E6" "5KuCoinTokenKCS"
Set-TextValue $ws "B7" "GateToken"
Set-TextValue $ws "C7" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D7" "3.045"
Set-TextValue $ws "E7" "6GateTokenGT"
Set-TextValue $ws "B8" "MXToken"
Set-TextValue $ws "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D8" "0.8106"
Set-TextValue $ws "E8" "7MXTokenMX"
Set-TextValue $ws "B9" "FTXToken"
Set-TextValue $ws "C9" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D9" "0.8381"
Set-TextValue $ws "E9" "8FTXTokenFTT"
Set-TextValue $ws "B10" "WazirX"
Set-TextValue $ws "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D10" "0.1334"
Set-TextValue $ws "E10" "9WazirXWRX"
Set-TextValue $ws "B11" "MandalaExchangeToken"
Set-TextValue $ws "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D11" "0.06956"
Set-TextValue $ws "E11" "10MandalaExchangeTokenMDX"
Set-TextValue $ws "B12" "BitrueCoin"
Set-TextValue $ws "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D12" "0.02815"
Set-TextValue $ws "E12" "11BitrueCoinBTR"
Set-TextValue $ws "B13" "BitMartToken"
Set-TextValue $ws "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D13" "0.09368"
Set-TextValue $ws "E13" "12BitMartTokenBMX"
Set-TextValue $ws "B14" "BitForexToken"
Set-TextValue $ws "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D14" "0.001521"
Set-TextValue $ws "E14" "13BitForexTokenBF"
Set-TextValue $ws "B15" "One"
Set-TextValue $ws "C15" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D15" "0.0005972"
Set-TextValue $ws "E15" "14OneONE"
Set-TextValue $ws "B16" "TigerCash"
Set-TextValue $ws "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D16" "0.006252"
Set-TextValue $ws "E16" "15TigerCashTCH"
Set-TextValue $ws "B17" "LEO"
Set-TextValue $ws "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D17" "3.499"
Set-TextValue $ws "E17" "16LEOLEO"
Set-TextValue $ws "D19" "0.3192"
Set-TextValue $ws "D20" "0.03221"
Set-TextValue $ws "D22" "3.757"
Set-TextValue $ws "D23" "0.04677"
Set-TextValue $ws "D25" "0.001234"
Set-TextValue $ws "D26" "0.004241"
Set-TextValue $ws "D27" "0.00009704"
Set-TextValue $ws "D40" "0.03622"
Set-TextValue $ws "B41" "BKEXToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D41" "0.1050"
Set-TextValue $ws "E41" "40BKEXTokenBKK"
Set-TextValue $ws "B42" "CEJI"
Set-TextValue $ws "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.002721"
Set-TextValue $ws "E42" "41CEJICEJIBestin24h"
Set-TextValue $ws "B43" "KickToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D43" "0.003233"
Set-TextValue $ws "E43" "42KickTokenKICK"
Set-TextValue $ws "D44" "0.007342"
Set-TextValue $ws "D45" "0.00005289"
Set-TextValue $ws "D49" "0.00002100"
Set-TextValue $ws "D50" "0.0002000"
